$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from
# 2023-09-13 (45182) to 2023-09-15 (45184) for every data row (2-253).
$startRow = 2
$endRow = 253

for ($row = $startRow; $row -le $endRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
